$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Exact "86÷7=" "92÷9="
Replace-Exact "83÷6=" "66÷8="
Replace-Exact "71÷5=" "72÷5="
Replace-Exact "32÷5=" "12÷8="
Replace-Exact "56÷4=" "70÷2="
Replace-Exact "61÷9=" "73÷4="
Replace-Exact "94÷8=" "27÷4="
Replace-Exact "63÷9=" "90÷7="
Replace-Exact "80÷5=" "72÷6="
Replace-Exact "17÷3=" "87÷3="
Replace-Exact "25÷7=" "80÷3="
Replace-Exact "34÷5=" "37÷4="
Replace-Exact "53÷7=" "71÷6="
Replace-Exact "93÷4=" "79÷6="
Replace-Exact "76÷6=" "63÷9="
Replace-Exact "75÷2=" "18÷8="
Replace-Exact "70÷9=" "73÷5="
Replace-Exact "26÷7=" "23÷5="
Replace-Exact "88÷9=" "18÷2="
Replace-Exact "99÷5=" "10÷2="
Replace-Exact "91÷3=" "53÷9="
Replace-Exact "40÷6=" "14÷6="
Replace-Exact "19÷2=" "88÷6="
Replace-Exact "89÷7=" "96÷8="
Replace-Exact "57÷4=" "61÷5="

Write-Output "Replacements complete"
